$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(''Kharis & The Beholder'', [''{1}{G}{G}{W}{W}'', ''Legendary Creature — Dragon Eye Wizard'', ''Flying'', ''When Kharis & The Beholder enters the battlefield and at the beginning of your upkeep, create a 1/1 white Human creature token and make a charisma check. (Roll a d20.)'', ''• If the result plus the number of creatures you control is greater than 11, put a +1/+1 counter on each creature you control.'', ''• If the result is a natural 20, for each nonlegendary creature you control, create a token that’s a copy of that creature.'', ''1/20''])'
$ws.Range("A3").Value = '(''Optimus Prime, Inspiring Leader'', [''{3}{R}{W}'', ''Autobot Character — Bot Mode'', ''{1}: Turn target permanent you control to its other face.'', ''{1}: Until end of turn, Optimus Prime, Inspiring Leader becomes a Construct with base power and toughness 6/6 and creatures you control gain trample.'', ''4/5''])'
$ws.Range("A4").Value = '(''Sol, Advocate Eternal'', [''{G}{W}{U}{B}'', ''Legendary Creature — Dragon Angel'', ''Legendary Partner (You can have two commanders if this is one of them. The other one is promoted to legendary.)'', ''Flying, vigilance'', ''Teamwork—Whenever you attack or block with both Sol, Advocate Eternal and its partner, support 4 and investigate four times.'', ''4/4''])'
$ws.Range("A5").Value = '(''The Legend of Arena'', [''{1}{U}{R}{W}'', ''Legendary Enchantment — Saga'', ''The Legend of Arena can be your commander.'', ''I, II — Create a 2/1 red Human Wizard creature token. Spells you cast this turn cost {1} less to cast for each Wizard you control.'', ''III — Search your library for a planeswalker card, put it onto the battlefield, then shuffle your library. It enters with an additional loyalty counter on it for each Wizard you control.''])'

$ws.Range("A6:A28").EntireRow.Delete()

Write-Output "Done. UsedRange: $($ws.UsedRange.Address())"
